$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(2, 39.74649, -74.20121),
    @(3, 39.78107, -104.91286),
    @(4, 35.47007, -82.87474),
    @(5, 46.8304, -100.77309),
    @(6, 43.81283, -94.07801000000001),
    @(7, 18.33551, -64.9636),
    @(8, 31.414179, -103.514546),
    @(9, 39.74874, -84.16482000000001),
    @(10, 26.065681, -80.23483400000001),
    @(11, 36.776494, -114.03794)
)

foreach ($row in $values) {
    $r = $row[0]
    $lat = $row[1]
    $lon = $row[2]
    $ws.Cells.Item($r, 1).Value = $lat
    $ws.Cells.Item($r, 2).Value = $lon
}
